# Fix column order for Birthdays sheet in Excel export
# New order: Date, Description, Grid Color, Text Color

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Birthdays")

# Header row
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Grid Color"
$ws.Range("D1").Value = "Text Color"

# Row 2 - Ahmet's Birthday
$ws.Range("B2").Value = "Ahmet's Birthday"
$ws.Range("C2").Value = "#87CEFA"
$ws.Range("D2").Value = "#000000"

# Row 3 - Buse's Birthday
$ws.Range("B3").Value = "Buse's Birthday"
$ws.Range("C3").Value = "#87CEFA"
$ws.Range("D3").Value = "#000000"

# Row 4 - Mehmet's Birthday
$ws.Range("B4").Value = "Mehmet's Birthday"
$ws.Range("C4").Value = "#87CEFA"
$ws.Range("D4").Value = "#000000"
